$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 2573.5
$ws.Range("J32").Value = 2573.5
$ws.Range("L32").Value = 2573.5
$ws.Range("N32").Value = -3225.5

# Row 38
$ws.Range("H38").Value = 694.4
$ws.Range("J38").Value = 988.625
$ws.Range("L38").Value = 2965.875
$ws.Range("N38").Value = -3709.875

# Row 129
$ws.Range("H129").Value = 862.8111
$ws.Range("J129").Value = 949.62823
$ws.Range("L129").Value = 2848.88469
$ws.Range("N129").Value = -12848.88469

# Row 137
$ws.Range("H137").Value = 1760.8823
$ws.Range("I137").Value = 1101.2
$ws.Range("K137").Value = 3303.6
$ws.Range("M137").Value = -753.6000000000004

# Row 138
$ws.Range("H138").Value = 1658.38
$ws.Range("I138").Value = 1208
$ws.Range("J138").Value = 1744.1666
$ws.Range("K138").Value = 3624
$ws.Range("L138").Value = 5232.4998
$ws.Range("M138").Value = 1516
$ws.Range("N138").Value = -15512.4998

# Row 139
$ws.Range("H139").Value = 34984
$ws.Range("J139").Value = 34984
$ws.Range("L139").Value = 34984
$ws.Range("N139").Value = -45264

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4469.7544
$ws.Range("J32").Value = 8722
$ws.Range("L32").Value = 8722
$ws.Range("N32").Value = -9296

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 4981.6875
$ws.Range("I86").Value = 5655.1816
$ws.Range("J86").Value = 3500
$ws.Range("K86").Value = 5655.1816
$ws.Range("L86").Value = 3500
$ws.Range("M86").Value = -4532.1816
$ws.Range("N86").Value = -5746

# Row 89
$ws.Range("H89").Value = 4981.6875
$ws.Range("I89").Value = 5655.1816
$ws.Range("J89").Value = 3500
$ws.Range("K89").Value = 28275.908
$ws.Range("L89").Value = 17500
$ws.Range("M89").Value = -22659.908
$ws.Range("N89").Value = -28732

# Row 99
$ws.Range("H99").Value = 41667770
$ws.Range("I99").Value = 50001020
$ws.Range("K99").Value = 50001020
$ws.Range("M99").Value = -49999522

# Row 134
$ws.Range("H134").Value = 4918.353
$ws.Range("I134").Value = 1203.6923
$ws.Range("J134").Value = 16991
$ws.Range("K134").Value = 3611.0769
$ws.Range("L134").Value = 50973
$ws.Range("M134").Value = -1076.0769
$ws.Range("N134").Value = -56043

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 83334584
$ws.Range("I16").Value = 111112310
$ws.Range("K16").Value = 111112310
$ws.Range("M16").Value = -111112023

# Row 31
$ws.Range("H31").Value = 1093.6428
$ws.Range("I31").Value = 1078.7273
$ws.Range("K31").Value = 1078.7273
$ws.Range("M31").Value = -783.7273

# Row 34
$ws.Range("H34").Value = 1093.6428
$ws.Range("I34").Value = 1078.7273
$ws.Range("K34").Value = 1078.7273
$ws.Range("M34").Value = -876.7273

# Row 58
$ws.Range("H58").Value = 719.46875
$ws.Range("I58").Value = 796.94446
$ws.Range("J58").Value = 619.8570999999999
$ws.Range("K58").Value = 796.94446
$ws.Range("L58").Value = 619.8570999999999
$ws.Range("M58").Value = -593.94446
$ws.Range("N58").Value = -1025.8571

# Row 86
$ws.Range("H86").Value = 2391006.5
$ws.Range("I86").Value = 3335776
$ws.Range("J86").Value = 29082.75
$ws.Range("K86").Value = 3335776
$ws.Range("L86").Value = 29082.75
$ws.Range("M86").Value = -3334653
$ws.Range("N86").Value = -31328.75

# Row 89
$ws.Range("H89").Value = 2391006.5
$ws.Range("I89").Value = 3335776
$ws.Range("J89").Value = 29082.75
$ws.Range("K89").Value = 16678880
$ws.Range("L89").Value = 145413.75
$ws.Range("M89").Value = -16673264
$ws.Range("N89").Value = -156645.75

# Row 113
$ws.Range("H113").Value = 83334584
$ws.Range("I113").Value = 111112310
$ws.Range("K113").Value = 111112310
$ws.Range("M113").Value = -111110140

# Row 132
$ws.Range("H132").Value = 1890.9166
$ws.Range("I132").Value = 1435.8334
$ws.Range("K132").Value = 4307.5002
$ws.Range("M132").Value = -1777.5002

# Row 134
$ws.Range("H134").Value = 796.23334
$ws.Range("I134").Value = 679.52
$ws.Range("J134").Value = 1379.8
$ws.Range("K134").Value = 2038.56
$ws.Range("L134").Value = 4139.4
$ws.Range("M134").Value = 496.4400000000001
$ws.Range("N134").Value = -9209.4

# Row 136
$ws.Range("H136").Value = 719.46875
$ws.Range("I136").Value = 796.94446
$ws.Range("J136").Value = 619.8570999999999
$ws.Range("K136").Value = 2390.83338
$ws.Range("L136").Value = 1859.5713
$ws.Range("M136").Value = 159.16662
$ws.Range("N136").Value = -6959.5713

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1334.8064
$ws.Range("J5").Value = 791.875
$ws.Range("L5").Value = 2375.625
$ws.Range("N5").Value = -2599.625

# Row 6
$ws.Range("H6").Value = 300.4
$ws.Range("I6").Value = 275
$ws.Range("J6").Value = 402
$ws.Range("K6").Value = 825
$ws.Range("L6").Value = 1206
$ws.Range("M6").Value = -712
$ws.Range("N6").Value = -1432

# Row 10
$ws.Range("H10").Value = 84.8
$ws.Range("I10").Value = 56
$ws.Range("J10").Value = 200
$ws.Range("K10").Value = 168
$ws.Range("L10").Value = 600
$ws.Range("M10").Value = -29
$ws.Range("N10").Value = -878

# Row 11
$ws.Range("H11").Value = 172288.28
$ws.Range("I11").Value = 193774.31
$ws.Range("K11").Value = 581322.9299999999
$ws.Range("M11").Value = -581182.9299999999

# Row 13
$ws.Range("H13").Value = 432.66666
$ws.Range("I13").Value = 149.5
$ws.Range("J13").Value = 999
$ws.Range("K13").Value = 448.5
$ws.Range("L13").Value = 2997
$ws.Range("M13").Value = -280.5
$ws.Range("N13").Value = -3333

# Row 135
$ws.Range("H135").Value = 1334.8064
$ws.Range("J135").Value = 791.875
$ws.Range("L135").Value = 7126.875
$ws.Range("N135").Value = -12196.875

# Row 140
$ws.Range("H140").Value = 21003.17
$ws.Range("J140").Value = 2903.5
$ws.Range("L140").Value = 8710.5
$ws.Range("N140").Value = -19070.5

$ws = $wb.Worksheets.Item("LTW")
# Row 64
$ws.Range("H64").Value = 19900
$ws.Range("J64").Value = 19900
$ws.Range("L64").Value = 19900
$ws.Range("N64").Value = -20350

# Row 67
$ws.Range("H67").Value = 19900
$ws.Range("J67").Value = 19900
$ws.Range("L67").Value = 19900
$ws.Range("N67").Value = -21460

# Row 136
$ws.Range("H136").Value = 1506.9048
$ws.Range("I136").Value = 1420.5
$ws.Range("J136").Value = 1622.1111
$ws.Range("K136").Value = 4261.5
$ws.Range("L136").Value = 4866.3333
$ws.Range("M136").Value = -1711.5
$ws.Range("N136").Value = -9966.3333

$ws = $wb.Worksheets.Item("WVR")
# Row 69
$ws.Range("H69").Value = 9450
$ws.Range("J69").Value = 9450
$ws.Range("L69").Value = 9450
$ws.Range("N69").Value = -10948

# Row 72
$ws.Range("H72").Value = 9450
$ws.Range("J72").Value = 9450
$ws.Range("L72").Value = 28350
$ws.Range("N72").Value = -35838

# Row 136
$ws.Range("H136").Value = 568.95
$ws.Range("I136").Value = 285.33334
$ws.Range("K136").Value = 856.0000200000001
$ws.Range("M136").Value = 1693.99998
